$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Yearly sheet: add a second "2017" block (columns J:O) mirroring the
# existing "2016" block (columns B:G).
# ---------------------------------------------------------------------------
$yearly = $wb.Worksheets.Item("Yearly")

# Copy the formatting of the existing B1:G15 block onto J1:O15 so every new
# cell picks up the same style indexes (fonts/fills/number formats) as its
# mirrored counterpart.
$yearly.Range("B1:G15").Copy()
$yearly.Range("J1").PasteSpecial(-4122)

# Year header
$yearly.Range("J1").Value = 2017

# Column headers (row 2)
$yearly.Range("K2").Value = "Month"
$yearly.Range("L2").Value = "Taxable Account"
$yearly.Range("M2").Value = "401K"
$yearly.Range("N2").Value = "Suzie's Roth IRA"
$yearly.Range("O2").Value = "Grand Total"

# Month labels + index numbers (rows 3-14)
$months = @("January","February","March","April","May","June","July","August","September","October","November","December")
for ($i = 0; $i -lt 12; $i++) {
    $r = 3 + $i
    $yearly.Range("J$r").Value = ($i + 1)
    $yearly.Range("K$r").Value = $months[$i]
}

# Monthly dollar amounts - only January (row 3) is non-zero
$yearly.Range("L3").Value = 24.7
$yearly.Range("M3").Value = 7.55
$yearly.Range("N3").Value = 0
for ($r = 4; $r -le 14; $r++) {
    $yearly.Range("L$r").Value = 0
    $yearly.Range("M$r").Value = 0
    $yearly.Range("N$r").Value = 0
}

# Grand total column (O) - row-wise sum, shared formula like column G
$yearly.Range("O3").Formula = "=SUM(L3:N3)"
for ($r = 4; $r -le 14; $r++) {
    $yearly.Range("O$r").Formula = "=SUM(L$r`:N$r)"
}

# Totals row (15)
$yearly.Range("K15").Value = "Total"
$yearly.Range("L15").Formula = "=SUM(L3:L14)"
$yearly.Range("M15").Formula = "=SUM(M3:M14)"
$yearly.Range("N15").Formula = "=SUM(N3:N14)"
$yearly.Range("O15").Formula = "=SUM(O3:O14)"

# ---------------------------------------------------------------------------
# All Time sheet: 2017's Taxable Account now pulls from Yearly!L3 instead of
# being a hard-coded literal.
# ---------------------------------------------------------------------------
$allTime = $wb.Worksheets.Item("All Time")
$allTime.Range("F8").Formula = "=Yearly!L3"

# ---------------------------------------------------------------------------
# View state: Yearly becomes the active/selected sheet (selection L4); All
# Time scrolls to show row 31 onward (selection L22) and is no longer active.
# ---------------------------------------------------------------------------
$yearly.Activate()
$yearly.Application.ActiveWindow.SelectedSheets.Item(1).Select()
$yearly.Range("L4").Select()

$allTime.Select()
$allTime.Application.ActiveWindow.ScrollRow = 31
$allTime.Range("L22").Select()

$yearly.Activate()
